# Updated symbol list (cryptos.xlsx) - refresh prices/links/volume labels
# and re-sort a few rows so KuCoinToken moves from row 23 up to row 7.
# Numeric-looking "Price" strings are written with a leading apostrophe
# (quote-prefix) so Excel keeps them as literal text instead of coercing
# them to numbers/doubles, then the style is reset back to Normal so no
# stray "Text" number-format / quote-prefix styling sticks to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'" + '242.97'
$ws.Range('D2').Style = 'Normal'
$ws.Range('D3').Value = "'" + '23.03'
$ws.Range('D3').Style = 'Normal'
$ws.Range('D4').Value = "'" + '5.414'
$ws.Range('D4').Style = 'Normal'
$ws.Range('D5').Value = "'" + '0.05924'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Value = "'" + '3.398'
$ws.Range('D6').Style = 'Normal'
$ws.Range('B7').Value = 'KuCoinToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range('D7').Value = "'" + '6.440'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '6KuCoinTokenKCS'
$ws.Range('B8').Value = 'MXToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D8').Value = "'" + '0.8042'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '7MXTokenMX'
$ws.Range('B9').Value = 'FTXToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D9').Value = "'" + '0.9117'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '8FTXTokenFTT'
$ws.Range('B10').Value = 'WazirX'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D10').Value = "'" + '0.1416'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '9WazirXWRX'
$ws.Range('B11').Value = 'MandalaExchangeToken'
$ws.Range('C11').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D11').Value = "'" + '0.07415'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '10MandalaExchangeTokenMDX'
$ws.Range('B12').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C12').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D12').Value = "'" + '0.03379'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '11LiechtensteinCryptoassetsExchangeLCX'
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D13').Value = "'" + '0.03088'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '12BitrueCoinBTR'
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D14').Value = "'" + '0.09327'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '13BitMartTokenBMX'
$ws.Range('B15').Value = 'MCDex'
$ws.Range('C15').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('D15').Value = "'" + '3.941'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '14MCDexMCB'
$ws.Range('B16').Value = 'BitForexToken'
$ws.Range('C16').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D16').Value = "'" + '0.001589'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '15BitForexTokenBF'
$ws.Range('B17').Value = 'CoinExToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D17').Value = "'" + '0.04808'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '16CoinExTokenCET'
$ws.Range('D18').Value = "'" + '0.005443'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').Value = "'" + '0.004288'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').Value = "'" + '0.0009837'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').Value = "'" + '0.00007521'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').Value = "'" + '3.651'
$ws.Range('D22').Style = 'Normal'
$ws.Range('B23').Value = 'BTSEToken'
$ws.Range('C23').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D23').Value = "'" + '2.185'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '22BTSETokenBTSE'
$ws.Range('B24').Value = 'One'
$ws.Range('C24').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D24').Value = "'" + '0.01112'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '23OneONEBestin24h'
$ws.Range('D25').Value = "'" + '0.3249'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').Value = "'" + '0.1349'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').Value = "'" + '0.0002452'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D40').Value = "'" + '0.03896'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').Value = "'" + '0.006233'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').Value = "'" + '0.1060'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').Value = "'" + '0.002749'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').Value = "'" + '0.006536'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '43LocalTradersLCT'
$ws.Range('D45').Value = "'" + '0.00005188'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').Value = "'" + '0.00000000752'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').Value = "'" + '0.0005811'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').Value = "'" + '1.052'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').Value = "'" + '0.002323'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').Value = "'" + '0.00002105'
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').Value = "'" + '0.0002005'
$ws.Range('D51').Style = 'Normal'
